$wb = $excel.ActiveWorkbook

# --- EURIBOR1M sheet: add helper formula in D2, update selection ---
$ws1 = $wb.Worksheets.Item("EURIBOR1M")
$ws1.Range("D2").Formula = "=B2*0.01"
$ws1.Range("D2").Style = "Normal"
$ws1.Range("E4").Select()

# --- EURIBOR6M sheet: it was the tab-selected sheet; no longer selected ---
$ws3 = $wb.Worksheets.Item("EURIBOR6M")
$ws3.Range("I17").Select()

# --- USDLIBOR3M sheet: multiply existing rates by 100, becomes the active tab ---
$ws4 = $wb.Worksheets.Item("USDLIBOR3M")
for ($r = 2; $r -le 18; $r++) {
    $cell = $ws4.Cells.Item($r, 2)
    $old = $cell.Value()
    $cell.Value = $old * 100
}
$ws4.Range("O15").Select()
$ws4.Activate()

$wb.Save()
